# Apply the "blurring / thresholding / label filtering" column additions
# plus the path-shortening edits to J4/L4, as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shorten the two existing long UNC-style paths in row 4 -----------------
$ws.Range("J4").Value = "...\dataset03\raw-cropNorm\"
$ws.Range("L4").Value = "...\dataset03\raw-cropNorm-bicubic-scaled0.25\"

# --- New headers (row 1), columns M:S ---------------------------------------
$ws.Range("M1").Value = "blurring sigma x"
$ws.Range("N1").Value = "blurring sigma y"
$ws.Range("O1").Value = "blurring sigma z"
$ws.Range("P1").Value = "blurring method"
$ws.Range("Q1").Value = "threshold label method"
$ws.Range("R1").Value = "threshold value (lower)"
$ws.Range("S1").Value = "label filtering/selection"

# --- Row 2 (dataset01) - values unknown, marked with "?" --------------------
$ws.Range("M2").Value = "?"
$ws.Range("N2").Value = "?"
$ws.Range("O2").Value = "?"
$ws.Range("P2").Value = "?"
$ws.Range("Q2").Value = "Otsu"
$ws.Range("R2").Value = "?"
$ws.Range("S2").Value = "?"

# --- Row 3 (dataset02) -------------------------------------------------------
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = "blur3D"
$ws.Range("Q3").Value = "Otsu"
$ws.Range("R3").Value = 570
$ws.Range("S3").Value = "Keep Largest Label"

# --- Row 4 (dataset03) -------------------------------------------------------
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = "blur3D"
$ws.Range("Q4").Value = "Otsu"
$ws.Range("R4").Value = 570
$ws.Range("S4").Value = "Keep Largest Label"

# --- Column width touch-ups to roughly mirror the authored layout -----------
$ws.Columns("I").ColumnWidth = 26.88671875
$ws.Columns("J").ColumnWidth = 25.6640625
$ws.Columns("L").ColumnWidth = 41.6640625
$ws.Columns("M").ColumnWidth = 13.6640625
$ws.Columns("N").ColumnWidth = 13.6640625
$ws.Columns("O").ColumnWidth = 13.6640625
$ws.Columns("P").ColumnWidth = 14
$ws.Columns("Q").ColumnWidth = 20
$ws.Columns("R").ColumnWidth = 19.88671875
$ws.Columns("S").ColumnWidth = 19.88671875

# --- Move the view / selection to mirror the saved state (topLeftCell M1) ---
try {
    $excel.ActiveWindow.ScrollColumn = 13
    $excel.ActiveWindow.ScrollRow = 1
} catch {}
$ws.Range("S3").Select()
